$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Append the new "Playtesting 2 / Presentation" time-card entries that
# were logged for 11/18/2014 - 12/2/2014 (rows 98-107).
# ---------------------------------------------------------------------

$rows = @(
  @{ Row=98;  Task="Trying to fix doors not using closed texture";                                   Hours=2;    Date=41961 },
  @{ Row=99;  Task="Playtesting 2";                                                                   Hours=3;    Date=41962 },
  @{ Row=100; Task="Fixing and improving from Playtest2";                                             Hours=7.25; Date=41963 },
  @{ Row=101; Task="Fixing and improving from Playtest2";                                             Hours=3;    Date=41964 },
  @{ Row=102; Task="Playtesting 2.1 - Playtesting with Jamies kid and getting feedback";              Hours=2;    Date=41964 },
  @{ Row=103; Task="Fixing and improving from Playtest2";                                             Hours=5;    Date=41966 },
  @{ Row=104; Task="Playtesting 2.2 - Playtesting with a couple of people";                           Hours=1;    Date=41967 },
  @{ Row=105; Task="Smooth Camera turn";                                                              Hours=2;    Date=41968 },
  @{ Row=106; Task="Presentation";                                                                    Hours=2;    Date=41969 }
)

foreach ($r in $rows) {
  $ws.Cells.Item($r.Row, 1).Value = $r.Task
  $ws.Cells.Item($r.Row, 2).Value = $r.Hours

  # Reuse the date number-format (style) that's already applied to the
  # date column above (C97) instead of minting a brand-new style.
  $ws.Cells.Item(97, 3).Copy()
  $ws.Cells.Item($r.Row, 3).PasteSpecial(-4122) | Out-Null
  $ws.Cells.Item($r.Row, 3).Value = $r.Date
}

# Final row: hours logged against a text date-range label rather than a
# single serial date.
$ws.Cells.Item(107, 1).Value = "Presentation+"
$ws.Cells.Item(107, 2).Value = 10
$ws.Cells.Item(107, 3).Value = "11/29/2014 - 12/2/2014"

$excel.CutCopyMode = 0

# A stray column D width tweak happened alongside the data entry.
$ws.Range("D:D").ColumnWidth = 12.95

# Move the selection/scroll position to reflect where the user ended up.
$ws.Range("C112").Select() | Out-Null
